$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update F2 formula to divide by 1000 (convert to meters)
$ws.Range("F2").Formula = "=E2*937/1000"

# New input parameters: Radius label + value
$ws.Range("H7").Value = "Radius"
$ws.Range("H8").Value = 0.0015874999999999999

# Update the selection to match the diff (H11)
$ws.Range("H11").Select()
